$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "[Small worlds](#sec:smallworlds)" link from A11 to A10
$ws.Range("A10").Value = "[Small worlds](#sec:smallworlds)"
$ws.Range("A11").Value = ""

# Remove the TODO placeholders in E10 and E11
$ws.Range("E10").Value = ""
$ws.Range("E11").Value = ""

# Add a link to the Hwk 4 homework document in G10
$ws.Range("G10").Value = "[Hwk 4: Problem set I](https://drive.google.com/file/d/1TUta8-8redraG0L044teOdA3SxX2eBtj/view?usp=sharing)"
